# Add data for 2021-10-04
# - Rename the "through September 25" sheet/header to "through September 26"
# - Update the counts in column B (the "through Sept NN" column) and a handful
#   of other month columns to reflect the one additional carjacking record
#   added for the date 2021-10-04 (which falls in neighborhoods touched below).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and the header label for the rolling "through" column.
$ws.Name = "Through 2021-09-26"
$ws.Range("B1").Value = "September 2021 (through September 26)"

# Garfield Park (row 2)
$ws.Range("B2").Value = 16
$ws.Range("BD2").Value = 3

# Humboldt Park (row 4)
$ws.Range("AU4").Value = 3

# Austin (row 5)
$ws.Range("B5").Value = 11

# Roseland (row 6)
$ws.Range("B6").Value = 6

# Little Village (row 8)
$ws.Range("K8").Value = 2

# Lower West Side (row 14)
$ws.Range("B14").Value = 2
$ws.Range("T14").Value = 1

# South Shore (row 17)
$ws.Range("B17").Value = 4

# Grand Boulevard (row 18)
$ws.Range("B18").Value = 3

# Wicker Park (row 19)
$ws.Range("B19").Value = 6

# Englewood (row 20)
$ws.Range("B20").Value = 4
$ws.Range("BD20").Value = 3

# River North (row 21)
$ws.Range("B21").Value = 2

# South Chicago (row 22)
$ws.Range("AC22").Value = 1

# Lincoln Park (row 30)
$ws.Range("AL30").Value = 1

# Hyde Park (row 34)
$ws.Range("B34").Value = 2

# Calumet Heights (row 36)
$ws.Range("AL36").Value = 1

# Washington Heights (row 41)
$ws.Range("T41").Value = 1
$ws.Range("AU41").Value = 2

# Streeterville (row 44)
$ws.Range("AL44").Value = 1

# Grand Crossing (row 55)
$ws.Range("B55").Value = 5

# Rogers Park (row 93)
$ws.Range("AU93").Value = 1
